# Update stats for 2025-12 (row 25 of the sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw figures for 2025-12
$schools = 6476
$authorities = 1007
$users = 6020604

$ws.Range("B25").Value = $schools
$ws.Range("C25").Value = $authorities
$ws.Range("D25").Value = $users

# Derived figures, recomputed from the raw inputs (row 13 = same month, prior year)
$schoolsPrior = $ws.Range("B13").Value2
$authoritiesPrior = $ws.Range("C13").Value2
$usersPrior = $ws.Range("D13").Value2

$usersPerSchool = $users / $schools
$yoySchools = ($schools - $schoolsPrior) / $schoolsPrior * 100
$yoyAuthorities = ($authorities - $authoritiesPrior) / $authoritiesPrior * 100
$yoyUsers = ($users - $usersPrior) / $usersPrior * 100

$ws.Range("E25").Value = $usersPerSchool
$ws.Range("F25").Value = $yoySchools
$ws.Range("G25").Value = $yoyAuthorities
$ws.Range("H25").Value = $yoyUsers
